$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H17").Value = 2060.3635
$ws_ALC.Range("I17").Value = 2000
$ws_ALC.Range("J17").Value = 2066.4
$ws_ALC.Range("K17").Value = 6000
$ws_ALC.Range("L17").Value = 6199.200000000001
$ws_ALC.Range("M17").Value = -5832
$ws_ALC.Range("N17").Value = -6535.200000000001
$ws_ALC.Range("H87").Value = 58999
$ws_ALC.Range("J87").Value = 58999
$ws_ALC.Range("L87").Value = 58999
$ws_ALC.Range("N87").Value = -61495
$ws_ALC.Range("H90").Value = 58999
$ws_ALC.Range("J90").Value = 58999
$ws_ALC.Range("L90").Value = 176997
$ws_ALC.Range("N90").Value = -189477
$ws_ALC.Range("H98").Value = 1110.4193
$ws_ALC.Range("I98").Value = 1044.5883
$ws_ALC.Range("J98").Value = 1190.3572
$ws_ALC.Range("K98").Value = 1044.5883
$ws_ALC.Range("L98").Value = 1190.3572
$ws_ALC.Range("M98").Value = 453.4117000000001
$ws_ALC.Range("N98").Value = -4186.3572
$ws_ALC.Range("H122").Value = 1110.4193
$ws_ALC.Range("I122").Value = 1044.5883
$ws_ALC.Range("J122").Value = 1190.3572
$ws_ALC.Range("K122").Value = 3133.7649
$ws_ALC.Range("L122").Value = 3571.0716
$ws_ALC.Range("M122").Value = -683.7648999999997
$ws_ALC.Range("N122").Value = -8471.071599999999
$ws_ALC.Range("H135").Value = 1243.2727
$ws_ALC.Range("I135").Value = 634.5
$ws_ALC.Range("K135").Value = 5710.5
$ws_ALC.Range("M135").Value = -3175.5
$ws_ALC.Range("H137").Value = 4588.0303
$ws_ALC.Range("I137").Value = 5767.826
$ws_ALC.Range("J137").Value = 1874.5
$ws_ALC.Range("K137").Value = 17303.478
$ws_ALC.Range("L137").Value = 5623.5
$ws_ALC.Range("M137").Value = -14753.478
$ws_ALC.Range("N137").Value = -10723.5
$ws_ALC.Range("H138").Value = 10992701
$ws_ALC.Range("I138").Value = 34485936
$ws_ALC.Range("J138").Value = 3929.5967
$ws_ALC.Range("K138").Value = 103457808
$ws_ALC.Range("L138").Value = 11788.7901
$ws_ALC.Range("M138").Value = -103452668
$ws_ALC.Range("N138").Value = -22068.7901
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 4354.76
$ws_ARM.Range("I32").Value = 3637.5454
$ws_ARM.Range("J32").Value = 9614.333000000001
$ws_ARM.Range("K32").Value = 3637.5454
$ws_ARM.Range("L32").Value = 9614.333000000001
$ws_ARM.Range("M32").Value = -3350.5454
$ws_ARM.Range("N32").Value = -10188.333
$ws_ARM.Range("H55").Value = 31199.6
$ws_ARM.Range("J55").Value = 37499.5
$ws_ARM.Range("L55").Value = 37499.5
$ws_ARM.Range("N55").Value = -38129.5
$ws_ARM.Range("H102").Value = 2017.25
$ws_ARM.Range("I102").Value = 1821.7
$ws_ARM.Range("J102").Value = 2995
$ws_ARM.Range("K102").Value = 1821.7
$ws_ARM.Range("L102").Value = 2995
$ws_ARM.Range("M102").Value = -199.7
$ws_ARM.Range("N102").Value = -6239
$ws_ARM.Range("H132").Value = 2241.7932
$ws_ARM.Range("I132").Value = 1806.4286
$ws_ARM.Range("K132").Value = 5419.2858
$ws_ARM.Range("M132").Value = -2889.2858
$ws_ARM.Range("H140").Value = 96259.75
$ws_ARM.Range("J140").Value = 97013.336
$ws_ARM.Range("L140").Value = 97013.336
$ws_ARM.Range("N140").Value = -107373.336
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H45").Value = 35000
$ws_BSM.Range("J45").Value = 35000
$ws_BSM.Range("L45").Value = 35000
$ws_BSM.Range("N45").Value = -36616
$ws_BSM.Range("H107").Value = 3337.5305
$ws_BSM.Range("I107").Value = 3095.3784
$ws_BSM.Range("J107").Value = 4084.1667
$ws_BSM.Range("K107").Value = 3095.3784
$ws_BSM.Range("L107").Value = 4084.1667
$ws_BSM.Range("M107").Value = -1175.3784
$ws_BSM.Range("N107").Value = -7924.1667
$ws_BSM.Range("H132").Value = 178448.67
$ws_BSM.Range("J132").Value = 178448.67
$ws_BSM.Range("L132").Value = 178448.67
$ws_BSM.Range("N132").Value = -188568.67
$ws_BSM.Range("H134").Value = 2475.8965
$ws_BSM.Range("I134").Value = 1849.08
$ws_BSM.Range("J134").Value = 6393.5
$ws_BSM.Range("K134").Value = 5547.24
$ws_BSM.Range("L134").Value = 19180.5
$ws_BSM.Range("M134").Value = -3012.24
$ws_BSM.Range("N134").Value = -24250.5
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 3940.8333
$ws_CRP.Range("J31").Value = 7500.4116
$ws_CRP.Range("L31").Value = 7500.4116
$ws_CRP.Range("N31").Value = -8090.4116
$ws_CRP.Range("H34").Value = 3940.8333
$ws_CRP.Range("J34").Value = 7500.4116
$ws_CRP.Range("L34").Value = 7500.4116
$ws_CRP.Range("N34").Value = -7904.4116
$ws_CRP.Range("H58").Value = 2366.35
$ws_CRP.Range("I58").Value = 1837.8572
$ws_CRP.Range("J58").Value = 3599.5
$ws_CRP.Range("K58").Value = 1837.8572
$ws_CRP.Range("L58").Value = 3599.5
$ws_CRP.Range("M58").Value = -1634.8572
$ws_CRP.Range("N58").Value = -4005.5
$ws_CRP.Range("H59").Value = 55000
$ws_CRP.Range("J59").Value = 55000
$ws_CRP.Range("L59").Value = 55000
$ws_CRP.Range("N59").Value = -57290
$ws_CRP.Range("H80").Value = 29332.334
$ws_CRP.Range("J80").Value = 29332.334
$ws_CRP.Range("L80").Value = 29332.334
$ws_CRP.Range("N80").Value = -31578.334
$ws_CRP.Range("H83").Value = 29332.334
$ws_CRP.Range("J83").Value = 29332.334
$ws_CRP.Range("L83").Value = 87997.00199999999
$ws_CRP.Range("N83").Value = -99229.00199999999
$ws_CRP.Range("H99").Value = 10804.454
$ws_CRP.Range("I99").Value = 7787.625
$ws_CRP.Range("J99").Value = 13643.823
$ws_CRP.Range("K99").Value = 7787.625
$ws_CRP.Range("L99").Value = 13643.823
$ws_CRP.Range("M99").Value = -6289.625
$ws_CRP.Range("N99").Value = -16639.823
$ws_CRP.Range("H126").Value = 10804.454
$ws_CRP.Range("I126").Value = 7787.625
$ws_CRP.Range("J126").Value = 13643.823
$ws_CRP.Range("K126").Value = 23362.875
$ws_CRP.Range("L126").Value = 40931.469
$ws_CRP.Range("M126").Value = -20892.875
$ws_CRP.Range("N126").Value = -45871.469
$ws_CRP.Range("H133").Value = 144108.95
$ws_CRP.Range("J133").Value = 144108.95
$ws_CRP.Range("L133").Value = 144108.95
$ws_CRP.Range("N133").Value = -149168.95
$ws_CRP.Range("H136").Value = 2366.35
$ws_CRP.Range("I136").Value = 1837.8572
$ws_CRP.Range("J136").Value = 3599.5
$ws_CRP.Range("K136").Value = 5513.571599999999
$ws_CRP.Range("L136").Value = 10798.5
$ws_CRP.Range("M136").Value = -2963.571599999999
$ws_CRP.Range("N136").Value = -15898.5
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H11").Value = 686.06665
$ws_CUL.Range("I11").Value = 401.2
$ws_CUL.Range("J11").Value = 1255.8
$ws_CUL.Range("K11").Value = 1203.6
$ws_CUL.Range("L11").Value = 3767.4
$ws_CUL.Range("M11").Value = -1063.6
$ws_CUL.Range("N11").Value = -4047.4
$ws_CUL.Range("H38").Value = 518.34485
$ws_CUL.Range("I38").Value = 380.4737
$ws_CUL.Range("J38").Value = 780.3
$ws_CUL.Range("K38").Value = 1141.4211
$ws_CUL.Range("L38").Value = 2340.9
$ws_CUL.Range("M38").Value = -794.4211
$ws_CUL.Range("N38").Value = -3034.9
$ws_CUL.Range("H50").Value = 2689
$ws_CUL.Range("I50").Value = 1599.6666
$ws_CUL.Range("J50").Value = 3097.5
$ws_CUL.Range("K50").Value = 4798.9998
$ws_CUL.Range("L50").Value = 9292.5
$ws_CUL.Range("M50").Value = -4317.9998
$ws_CUL.Range("N50").Value = -10254.5
$ws_CUL.Range("H53").Value = 2689
$ws_CUL.Range("I53").Value = 1599.6666
$ws_CUL.Range("J53").Value = 3097.5
$ws_CUL.Range("K53").Value = 4798.9998
$ws_CUL.Range("L53").Value = 9292.5
$ws_CUL.Range("M53").Value = -4317.9998
$ws_CUL.Range("N53").Value = -10254.5
$ws_CUL.Range("H57").Value = 11375.647
$ws_CUL.Range("I57").Value = 8298.25
$ws_CUL.Range("K57").Value = 24894.75
$ws_CUL.Range("M57").Value = -24335.75
$ws_CUL.Range("H94").Value = 3000000
$ws_CUL.Range("I94").Value = 0
$ws_CUL.Range("K94").Value = 0
$ws_CUL.Range("M94").ClearContents()
$ws_CUL.Range("H122").Value = 839.5714
$ws_CUL.Range("J122").Value = 897.5
$ws_CUL.Range("L122").Value = 8077.5
$ws_CUL.Range("N122").Value = -12977.5
$ws_CUL.Range("H133").Value = 0
$ws_CUL.Range("I133").Value = 0
$ws_CUL.Range("K133").Value = 0
$ws_CUL.Range("M133").ClearContents()
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H62").Value = 44998.5
$ws_GSM.Range("J62").Value = 44998.5
$ws_GSM.Range("L62").Value = 44998.5
$ws_GSM.Range("N62").Value = -46370.5
$ws_GSM.Range("H65").Value = 44998.5
$ws_GSM.Range("J65").Value = 44998.5
$ws_GSM.Range("L65").Value = 134995.5
$ws_GSM.Range("N65").Value = -141859.5
$ws_GSM.Range("H97").Value = 540.1
$ws_GSM.Range("I97").Value = 496.36365
$ws_GSM.Range("J97").Value = 660.375
$ws_GSM.Range("K97").Value = 496.36365
$ws_GSM.Range("L97").Value = 660.375
$ws_GSM.Range("M97").Value = -0.3636500000000069
$ws_GSM.Range("N97").Value = -1652.375
$ws_GSM.Range("H122").Value = 2416.5454
$ws_GSM.Range("I122").Value = 1821.6
$ws_GSM.Range("K122").Value = 5464.799999999999
$ws_GSM.Range("M122").Value = -3014.799999999999
$ws_GSM.Range("H126").Value = 8413.223
$ws_GSM.Range("I126").Value = 8277
$ws_GSM.Range("J126").Value = 8627.286
$ws_GSM.Range("K126").Value = 24831
$ws_GSM.Range("L126").Value = 25881.858
$ws_GSM.Range("M126").Value = -22361
$ws_GSM.Range("N126").Value = -30821.858
$ws_GSM.Range("H132").Value = 2456.9153
$ws_GSM.Range("I132").Value = 2100.7708
$ws_GSM.Range("K132").Value = 6302.312399999999
$ws_GSM.Range("M132").Value = -3772.312399999999
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H99").Value = 43085.375
$ws_LTW.Range("I99").Value = 30997
$ws_LTW.Range("K99").Value = 30997
$ws_LTW.Range("M99").Value = -28002
$ws_LTW.Range("H133").Value = 49987.168
$ws_LTW.Range("J133").Value = 49987.168
$ws_LTW.Range("L133").Value = 49987.168
$ws_LTW.Range("N133").Value = -55047.168
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H62").Value = 23264110
$ws_WVR.Range("I62").Value = 24398176
$ws_WVR.Range("J62").Value = 15749
$ws_WVR.Range("K62").Value = 24398176
$ws_WVR.Range("L62").Value = 15749
$ws_WVR.Range("M62").Value = -24397552
$ws_WVR.Range("N62").Value = -16997
$ws_WVR.Range("H65").Value = 23264110
$ws_WVR.Range("I65").Value = 24398176
$ws_WVR.Range("J65").Value = 15749
$ws_WVR.Range("K65").Value = 121990880
$ws_WVR.Range("L65").Value = 78745
$ws_WVR.Range("M65").Value = -121987760
$ws_WVR.Range("N65").Value = -84985
$ws_WVR.Range("H123").Value = 198460.67
$ws_WVR.Range("J123").Value = 198460.67
$ws_WVR.Range("L123").Value = 198460.67
$ws_WVR.Range("N123").Value = -208260.67
$ws_WVR.Range("H132").Value = 3034.739
$ws_WVR.Range("I132").Value = 3091.9
$ws_WVR.Range("J132").Value = 2653.6667
$ws_WVR.Range("K132").Value = 9275.700000000001
$ws_WVR.Range("L132").Value = 7961.000100000001
$ws_WVR.Range("M132").Value = -6745.700000000001
$ws_WVR.Range("N132").Value = -13021.0001
